$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.081.38"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "2.631.34"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'310.34"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").Value = "'98.67"
$ws.Range("E6").Value = "  -4.80%  "
$ws.Range("E7").Value = "  -1.99%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "'38.75"
$ws.Range("E10").Value = "  -1.41%  "
$ws.Range("D11").Value = "'0.0846"
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'54.25"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "'8.07"
$ws.Range("E13").Value = "  -3.60%  "
$ws.Range("D14").Value = "3.028.64"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").Value = "2.631.62"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "46.046.55"
$ws.Range("E19").Value = "  -3.05%  "
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").Value = "'12.76"
$ws.Range("D23").Value = "'74.85"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("D24").Value = "'282.85"
$ws.Range("E24").Value = "  +8.35%  "
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'30.00"
$ws.Range("E27").Value = "  -4.56%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").Value = "'10.56"
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "'38.64"
$ws.Range("E30").Value = "  -6.66%  "
$ws.Range("E31").Value = "  -4.89%  "
$ws.Range("D32").Value = "'6.25"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("D33").Value = "'3.71"
$ws.Range("E33").Value = "  -2.52%  "
$ws.Range("D34").Value = "'2.33"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "'156.70"
$ws.Range("E35").Value = "  +2.19%  "
$ws.Range("D36").Value = "'0.0841"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -2.22%  "
$ws.Range("D38").Value = "'0.124"
$ws.Range("E38").Value = "  +4.09%  "
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").Value = "'22.47"
$ws.Range("E40").Value = "  +3.42%  "
$ws.Range("D41").Value = "'15.79"
$ws.Range("E41").Value = "  -7.89%  "
$ws.Range("D42").Value = "'0.0329"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'3.57"
$ws.Range("E43").Value = "  -4.19%  "
$ws.Range("D44").Value = "'4.04"
$ws.Range("E44").Value = "  -7.57%  "
$ws.Range("D45").Value = "2.126.69"
$ws.Range("E45").Value = "  +3.13%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'94.05"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'110.32"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").Value = "'9.14"
$ws.Range("E49").Value = "  -1.86%  "
$ws.Range("D50").Value = "2.879.92"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").Value = "'0.201"
$ws.Range("E51").Value = "  -1.31%  "
